# Auto-generated Excel COM-interop edit script
# Commit: Update automatic: dades i banners [2026-03-01 02:50]
# Updates DATA_EXTRACCIO timestamps and re-extracted meteo readings
# for rows 2-46 of the "Dades_Meteo" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / numeric-looking-but-not-pure-percent cells ---
# Assigning these directly as .Value is safe: they contain letters,
# degree signs, or other non-numeric characters so Excel keeps them
# as literal text (matches original inlineStr storage).
$ws.Range("E2").Value = "2026-03-01 02:48:32"
$ws.Range("O2").Value = "-0.4 °C"
$ws.Range("E3").Value = "2026-03-01 02:48:35"
$ws.Range("O3").Value = "-3.4 °C"
$ws.Range("E4").Value = "2026-03-01 02:48:37"
$ws.Range("M4").Value = "8.6 °C 2:28 TU"
$ws.Range("O4").Value = "8.3 °C"
$ws.Range("E5").Value = "2026-03-01 02:48:39"
$ws.Range("N5").Value = "-3.6 °C 2:28 TU"
$ws.Range("O5").Value = "-3.3 °C"
$ws.Range("E6").Value = "2026-03-01 02:48:42"
$ws.Range("J6").Value = "1025.6 hPa"
$ws.Range("O6").Value = "9.5 °C"
$ws.Range("E7").Value = "2026-03-01 02:48:44"
$ws.Range("E8").Value = "2026-03-01 02:48:46"
$ws.Range("J8").Value = "1025.7 hPa"
$ws.Range("O8").Value = "9.5 °C"
$ws.Range("E9").Value = "2026-03-01 02:48:49"
$ws.Range("E10").Value = "2026-03-01 02:48:52"
$ws.Range("M10").Value = "7.8 °C 2:26 TU"
$ws.Range("O10").Value = "6.6 °C"
$ws.Range("E11").Value = "2026-03-01 02:48:54"
$ws.Range("N11").Value = "6.1 °C 2:02 TU"
$ws.Range("O11").Value = "6.3 °C"
$ws.Range("E12").Value = "2026-03-01 02:48:57"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-03-01 02:48:59"
$ws.Range("N13").Value = "4.3 °C 2:29 TU"
$ws.Range("O13").Value = "4.5 °C"
$ws.Range("E14").Value = "2026-03-01 02:49:01"
$ws.Range("O14").Value = "11.4 °C"
$ws.Range("E15").Value = "2026-03-01 02:49:04"
$ws.Range("O15").Value = "7.9 °C"
$ws.Range("E16").Value = "2026-03-01 02:49:06"
$ws.Range("E17").Value = "2026-03-01 02:49:09"
$ws.Range("N17").Value = "1.4 °C 2:17 TU"
$ws.Range("E18").Value = "2026-03-01 02:49:12"
$ws.Range("J18").Value = "1025.9 hPa"
$ws.Range("E19").Value = "2026-03-01 02:49:14"
$ws.Range("N19").Value = "6.0 °C 2:00 TU"
$ws.Range("E20").Value = "2026-03-01 02:49:16"
$ws.Range("E21").Value = "2026-03-01 02:49:18"
$ws.Range("N21").Value = "6.3 °C 2:11 TU"
$ws.Range("O21").Value = "6.9 °C"
$ws.Range("E22").Value = "2026-03-01 02:49:21"
$ws.Range("N22").Value = "-5.4 °C 2:29 TU"
$ws.Range("O22").Value = "-4.8 °C"
$ws.Range("E23").Value = "2026-03-01 02:49:24"
$ws.Range("E24").Value = "2026-03-01 02:49:26"
$ws.Range("J24").Value = "1026.7 hPa"
$ws.Range("O24").Value = "4.5 °C"
$ws.Range("E25").Value = "2026-03-01 02:49:29"
$ws.Range("L25").Value = "9.4 km/h - 238º 2:07 TU"
$ws.Range("N25").Value = "-2.0 °C 2:28 TU"
$ws.Range("E26").Value = "2026-03-01 02:49:32"
$ws.Range("N26").Value = "2.5 °C 2:11 TU"
$ws.Range("O26").Value = "2.7 °C"
$ws.Range("E27").Value = "2026-03-01 02:49:34"
$ws.Range("N27").Value = "-1.3 °C 2:23 TU"
$ws.Range("E28").Value = "2026-03-01 02:49:37"
$ws.Range("J28").Value = "1025.7 hPa"
$ws.Range("N28").Value = "8.4 °C 2:26 TU"
$ws.Range("E29").Value = "2026-03-01 02:49:39"
$ws.Range("N29").Value = "9.1 °C 2:03 TU"
$ws.Range("O29").Value = "9.6 °C"
$ws.Range("E30").Value = "2026-03-01 02:49:42"
$ws.Range("J30").Value = "1025.7 hPa"
$ws.Range("O30").Value = "9.5 °C"
$ws.Range("E31").Value = "2026-03-01 02:49:44"
$ws.Range("E32").Value = "2026-03-01 02:49:47"
$ws.Range("O32").Value = "1.5 °C"
$ws.Range("E33").Value = "2026-03-01 02:49:50"
$ws.Range("J33").Value = "1026.0 hPa"
$ws.Range("M33").Value = "4.1 °C 2:10 TU"
$ws.Range("E34").Value = "2026-03-01 02:49:52"
$ws.Range("E35").Value = "2026-03-01 02:49:55"
$ws.Range("E36").Value = "2026-03-01 02:49:58"
$ws.Range("O36").Value = "9.1 °C"
$ws.Range("E37").Value = "2026-03-01 02:50:00"
$ws.Range("J37").Value = "1026.8 hPa"
$ws.Range("N37").Value = "6.2 °C 2:29 TU"
$ws.Range("E38").Value = "2026-03-01 02:50:03"
$ws.Range("M38").Value = "9.2 °C 2:08 TU"
$ws.Range("O38").Value = "8.5 °C"
$ws.Range("E39").Value = "2026-03-01 02:50:06"
$ws.Range("E40").Value = "2026-03-01 02:50:08"
$ws.Range("N40").Value = "6.9 °C 2:01 TU"
$ws.Range("O40").Value = "7.5 °C"
$ws.Range("E41").Value = "2026-03-01 02:50:10"
$ws.Range("E42").Value = "2026-03-01 02:50:12"
$ws.Range("O42").Value = "10.2 °C"
$ws.Range("E43").Value = "2026-03-01 02:50:15"
$ws.Range("N43").Value = "8.6 °C 2:18 TU"
$ws.Range("O43").Value = "8.8 °C"
$ws.Range("E44").Value = "2026-03-01 02:50:17"
$ws.Range("N44").Value = "-2.7 °C 2:20 TU"
$ws.Range("O44").Value = "-2.2 °C"
$ws.Range("E45").Value = "2026-03-01 02:50:20"
$ws.Range("J45").Value = "1027.2 hPa"
$ws.Range("N45").Value = "3.7 °C 2:28 TU"
$ws.Range("E46").Value = "2026-03-01 02:50:23"
$ws.Range("O46").Value = "7.3 °C"

# --- Pure percentage cells ("92%", "73%", ...) ---
# Assigning a bare "NN%" string straight to .Value makes Excel
# auto-convert it to a numeric percentage (changes cell type AND
# style/number-format), which does not match the source workbook
# where these remain literal text. Instead we build the literal
# string via a concatenation formula (so Excel never number-sniffs
# it), then Copy + PasteSpecial(values) to collapse the formula
# down to a static value while preserving the original cell style.
$ws.Range("H2").Formula = '="92" & "%"'
$ws.Range("H2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4163) | Out-Null
$ws.Range("H7").Formula = '="73" & "%"'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4163) | Out-Null
$ws.Range("H8").Formula = '="91" & "%"'
$ws.Range("H8").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4163) | Out-Null
$ws.Range("H9").Formula = '="63" & "%"'
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4163) | Out-Null
$ws.Range("H12").Formula = '="71" & "%"'
$ws.Range("H12").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4163) | Out-Null
$ws.Range("H14").Formula = '="96" & "%"'
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4163) | Out-Null
$ws.Range("H16").Formula = '="88" & "%"'
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4163) | Out-Null
$ws.Range("H23").Formula = '="94" & "%"'
$ws.Range("H23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4163) | Out-Null
$ws.Range("H26").Formula = '="96" & "%"'
$ws.Range("H26").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4163) | Out-Null
$ws.Range("H29").Formula = '="78" & "%"'
$ws.Range("H29").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4163) | Out-Null
$ws.Range("H30").Formula = '="88" & "%"'
$ws.Range("H30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4163) | Out-Null
$ws.Range("H36").Formula = '="81" & "%"'
$ws.Range("H36").Copy() | Out-Null
$ws.Range("H36").PasteSpecial(-4163) | Out-Null
$ws.Range("H40").Formula = '="85" & "%"'
$ws.Range("H40").Copy() | Out-Null
$ws.Range("H40").PasteSpecial(-4163) | Out-Null
$ws.Range("H41").Formula = '="87" & "%"'
$ws.Range("H41").Copy() | Out-Null
$ws.Range("H41").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
